$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 246, pushing the existing row 246 (and
# everything below it) down by one. This mirrors the diff: a brand new
# weekly record is inserted, and all subsequent rows (previously 246-256)
# shift to 247-257.
$ws.Rows.Item(246).Insert()

# Populate the newly inserted row 246 with the new weekly record. The
# "metadata" columns (market/region/category/etc.) are identical to the
# row that used to occupy this slot, only the date and price figures
# (D, J, K, L, M, P) are new.
$ws.Cells.Item(246, 1).Value = 10
$ws.Cells.Item(246, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(246, 3).Value = "La Araucanía"
$ws.Cells.Item(246, 4).Value = 44747
$ws.Cells.Item(246, 5).Value = 9
$ws.Cells.Item(246, 6).Value = 100112039
$ws.Cells.Item(246, 7).Value = "Ciboulette"
$ws.Cells.Item(246, 8).Value = "Sin especificar"
$ws.Cells.Item(246, 9).Value = "Primera"
$ws.Cells.Item(246, 10).Value = 55
$ws.Cells.Item(246, 11).Value = 9000
$ws.Cells.Item(246, 12).Value = 10000
$ws.Cells.Item(246, 13).Value = 9636
$ws.Cells.Item(246, 14).Value = "$/docena de atados"
$ws.Cells.Item(246, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(246, 16).Value = 3212
$ws.Cells.Item(246, 17).Value = 3
$ws.Cells.Item(246, 18).Value = "Hortaliza"
